$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.837.68"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "2.499.09"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("D9").Value = "2.499.11"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("D15").Value = "2.954.78"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "67.636.26"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "2.494.96"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -6.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "2.608.42"
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.549"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").Value = "0.0₆0275"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("E51").Value = "  -1.92%  "
